$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "URL: https://produto.mercadolivre.com.br/MLB-3519051321-fonte-carregador-automotivo-jfa-120a-bob-bivolt-automatico-_JM",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "URL: https://produto.mercadolivre.com.br/MLB-3753009491-fonte-carregador-jfa-120a-bob-slim-bivolt-cor-preto-bob120-_JM?searchVariation=183276947923",
    2)

$d.Content.Find.Execute(
    "Nome: Fonte Carregador Automotivo Jfa 120a Bob Bivolt Automático",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Nome: Fonte Carregador Jfa 120a Bob Slim Bivolt Cor Preto Bob120",
    2)

$d.Content.Find.Execute(
    "Preço: 514.44",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Preço: 554.96",
    2)

$d.Content.Find.Execute(
    "Preço Previsto: 514.45",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Preço Previsto: 555.93",
    2)

$d.Content.Find.Execute(
    "Loja: Motor Shop",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Loja: FLORIPASOUND",
    2)

$d.Content.Find.Execute(
    "Tipo: Clássico",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Tipo: Premium",
    2)

$d.Content.Find.Execute(
    "Lugar: Pedra branca, Ceará.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Lugar: Palhoça, Santa Catarina.",
    2)
